$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Icon" column (C) with the new product image file names.
# Mapping determined by category name (column B) -> new icon file name.
$ws.Range("C2").Value  = "icon-11.png"  # Buah-Buahan
$ws.Range("C3").Value  = "icon-12.png"  # Sayur Segar
$ws.Range("C4").Value  = "icon-10.png"  # Jajanan Pasar
$ws.Range("C5").Value  = "icon-9.png"   # Bumbu Dapur & Kue
$ws.Range("C6").Value  = "icon-13.png"  # Ayam
$ws.Range("C7").Value  = "icon-14.png"  # Ikan & Seafood
$ws.Range("C8").Value  = "icon-3.png"   # Daging
$ws.Range("C9").Value  = "icon-1.png"   # Tahu, Tempe dan Nabati
$ws.Range("C10").Value = "icon-8.png"   # Siap Saji dan Olahan
$ws.Range("C11").Value = "icon-6.png"   # Kopi, Teh dan Minuman
$ws.Range("C12").Value = "icon-8.png"   # Susu dan Olahan Susu
$ws.Range("C13").Value = "icon-4.png"   # Rumah Tangga
$ws.Range("C14").Value = "icon-5.png"   # Kesehatan
$ws.Range("C15").Value = "icon-2.png"   # Perlengkapan Makanan

# Move the active selection to C14, matching the saved cursor position.
$ws.Range("C14").Select()
